# Added 4wk low sales check: updates MyForecast (D), Inventory Coverage (H),
# Stockout Risk (I), and Seasonality Index (L) values on the
# "Forecast Comparison" sheet, plus the derived aggregate metrics on the
# "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet -------------------------------------------
# Row 2 (W10)
$wsForecast.Range("D2").Value = 44
$wsForecast.Range("H2").Value = 2.91
$wsForecast.Range("L2").Value = 1.15

# Row 3 (W11)
$wsForecast.Range("D3").Value = 41
$wsForecast.Range("H3").Value = 2.05
$wsForecast.Range("L3").Value = 1.16

# Row 4 (W12)
$wsForecast.Range("D4").Value = 37
$wsForecast.Range("H4").Value = 1.16
$wsForecast.Range("L4").Value = 1.06

# Row 5 (W13)
$wsForecast.Range("D5").Value = 35
$wsForecast.Range("H5").Value = 0.17
$wsForecast.Range("I5").Value = "High"
$wsForecast.Range("L5").Value = 1.18

# Row 6 (W14)
$wsForecast.Range("D6").Value = 36
$wsForecast.Range("L6").Value = 1.1

# Row 7 (W15)
$wsForecast.Range("D7").Value = 34
$wsForecast.Range("L7").Value = 0.82

# Row 8 (W16)
$wsForecast.Range("D8").Value = 33
$wsForecast.Range("L8").Value = 0.98

# Row 9 (W17)
$wsForecast.Range("L9").Value = 0.89

# Row 10 (W18)
$wsForecast.Range("D10").Value = 34
$wsForecast.Range("L10").Value = 0.96

# Row 11 (W19)
$wsForecast.Range("L11").Value = 1.09

# Row 12 (W20)
$wsForecast.Range("D12").Value = 34
$wsForecast.Range("L12").Value = 1.01

# Row 13 (W21)
$wsForecast.Range("D13").Value = 33
$wsForecast.Range("L13").Value = 1.19

# Row 14 (W22)
$wsForecast.Range("D14").Value = 32
$wsForecast.Range("L14").Value = 1.08

# Row 15 (W23)
$wsForecast.Range("D15").Value = 28
$wsForecast.Range("L15").Value = 0.86

# Row 16 (W24)
$wsForecast.Range("D16").Value = 28
$wsForecast.Range("L16").Value = 1.13

# Row 17 (W25)
$wsForecast.Range("D17").Value = 25
$wsForecast.Range("L17").Value = 1.09

# --- Summary sheet ---------------------------------------------------------
# These cells store their numbers as text (inlineStr) in the workbook, so
# force a text number format before assigning the value - otherwise Excel
# would auto-coerce the numeric-looking string into a real number.
$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "539"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "293"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "157"

$wsSummary.Range("B12").NumberFormat = "@"
$wsSummary.Range("B12").Value = "44"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "25"
